# Update "想去人数" (F column) counts that changed between scraped data runs.
# Commit: "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value  = 6968
$wsExhibition.Range("F6").Value  = 553
$wsExhibition.Range("F9").Value  = 4619
$wsExhibition.Range("F13").Value = 1413
$wsExhibition.Range("F15").Value = 125
$wsExhibition.Range("F32").Value = 113
$wsExhibition.Range("F42").Value = 326
$wsExhibition.Range("F43").Value = 1186

# Sheet "本地生活" (sheet3)
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F6").Value = 619
$wsLocal.Range("F8").Value = 1357
$wsLocal.Range("F9").Value = 2136

# Sheet "全部类型" (sheet4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value  = 619
$wsAll.Range("F9").Value  = 619
$wsAll.Range("F11").Value = 6968
$wsAll.Range("F12").Value = 553
$wsAll.Range("F14").Value = 4619
$wsAll.Range("F18").Value = 1413
$wsAll.Range("F21").Value = 125
$wsAll.Range("F22").Value = 1357
$wsAll.Range("F23").Value = 2136
$wsAll.Range("F35").Value = 113
$wsAll.Range("F45").Value = 326
